$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.483.51"
$ws.Range("E2").Value = "  +4.54%  "
$ws.Range("D3").Value = "3.627.36"
$ws.Range("E3").Value = "  +7.00%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").Value = "'183.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("D7").Value = "3.618.25"
$ws.Range("E7").Value = "  +7.04%  "
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("E10").Value = "  +4.94%  "
$ws.Range("E11").Value = "  +2.49%  "
$ws.Range("D12").Value = "'50.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.06%  "
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").Value = "'709.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.21%  "
$ws.Range("D15").Value = "4.208.21"
$ws.Range("E15").Value = "  +7.04%  "
$ws.Range("E16").Value = "  +3.60%  "
$ws.Range("D17").Value = "72.514.84"
$ws.Range("E17").Value = "  +4.58%  "
$ws.Range("D18").Value = "3.620.81"
$ws.Range("E18").Value = "  +6.55%  "
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").Value = "'18.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.78%  "
$ws.Range("E21").Value = "  +3.59%  "
$ws.Range("E22").Value = "  +3.18%  "
$ws.Range("E23").Value = "  +9.22%  "
$ws.Range("D24").Value = "'17.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.19%  "
$ws.Range("D25").Value = "'105.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.16%  "
$ws.Range("E26").Value = "  +2.95%  "
$ws.Range("E27").Value = "  +4.91%  "
$ws.Range("D28").Value = "'10.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.33%  "
$ws.Range("D29").Value = "'35.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.75%  "
$ws.Range("E30").Value = "  +4.50%  "
$ws.Range("E31").Value = "  +6.68%  "
$ws.Range("D32").Value = "'4.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +14.09%  "
$ws.Range("D33").Value = "'595.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.09%  "
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("D36").Value = "'59.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.145"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.64%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "3.644.73"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("E40").Value = "  +8.58%  "
$ws.Range("D41").Value = "'35.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  +6.78%  "
$ws.Range("E43").Value = "  +4.84%  "
$ws.Range("E44").Value = "  +6.52%  "
$ws.Range("D45").Value = "'0.353"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.91%  "
$ws.Range("D46").Value = "'3.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("E47").Value = "  +5.25%  "
$ws.Range("D48").Value = "'1.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.56%  "
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "'133.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.11%  "
